# Reports and Driver Page Automation Script Added
#
# Updates the sample Drivers record (Cell + License number) and
# re-selects B5 on the Drivers sheet, and normalizes the duplicate
# "Text" number-format style used on the Company_Profile sheet so it
# shares the same style slot as the rest of the workbook.

$wb = $excel.ActiveWorkbook

# --- Drivers sheet -------------------------------------------------
$drivers = $wb.Worksheets.Item("Drivers")

# Cell number
$drivers.Range("B3").Value = "911-000-000"

# License number
$drivers.Range("B5").Value = "AI1111111111"

# Leave the sheet's selection on B5 (matches the saved workbook state)
$drivers.Range("B5").Select()

# --- Company_Profile sheet ------------------------------------------
# A1:B3 used a duplicate "Text" (@) number-format style; re-apply the
# text format so it collapses onto the single shared style.
$profile = $wb.Worksheets.Item("Company_Profile")
$profile.Range("A1:B3").NumberFormat = "@"
